# Generate Report for Handoff
#
# The localization-status report tracks, per target file, the datetime of
# its latest handoff. Row 5 corresponds to the
# "8870a6a3-3da5-4b58-a93b-f0e6c6ce4918" source file, which has just been
# handed off for translation in both the zh-cn and de-de sheets. Stamp the
# "Latest Handoff Datetime" column (D) with the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-23 03:53:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-23 03:54:01"
